# Update Dll1-Notch2 sheet values with recalculated TPM-derived statistics
$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("G2").Value = 11.29255533333333
$ws.Range("H2").Value = 33.877666
$ws.Range("I2").Value = 0.5495662219753726
$ws.Range("J2").Value = 0.6375557499803809
$ws.Range("M2").Value = 1.400501333333333
$ws.Range("N2").Value = 4.201504
$ws.Range("O2").Value = 0.00926314904242919
$ws.Range("P2").Value = 0.009687730200823723
$ws.Range("Q2").Value = 15.81523880107378
$ws.Range("R2").Value = 142.337149209664
$ws.Range("S2").Value = 0.005090713822842601
$ws.Range("T2").Value = 0.006176468093793754

$ws.Range("G3").Value = 11.29255533333333
$ws.Range("H3").Value = 33.877666
$ws.Range("I3").Value = 0.5495662219753726
$ws.Range("J3").Value = 0.6375557499803809
$ws.Range("O3").Value = 0.1405812059498714
$ws.Range("P3").Value = 0.1470248171880475
$ws.Range("Q3").Value = 240.0183061782087
$ws.Range("R3").Value = 2160.164755603878
$ws.Range("S3").Value = 0.0772586822346126
$ws.Range("T3").Value = 0.093736517588054

$ws.Range("G4").Value = 11.29255533333333
$ws.Range("H4").Value = 33.877666
$ws.Range("I4").Value = 0.5495662219753726
$ws.Range("J4").Value = 0.6375557499803809
$ws.Range("M4").Value = 63.87756733333333
$ws.Range("N4").Value = 191.632702
$ws.Range("O4").Value = 0.4224968677952986
$ws.Range("P4").Value = 0.4418622271050682
$ws.Range("Q4").Value = 721.3409636703925
$ws.Range("R4").Value = 6492.068673033533
$ws.Range("S4").Value = 0.2321900074306907
$ws.Range("T4").Value = 0.2817118035899732

$ws.Range("G5").Value = 11.29255533333333
$ws.Range("H5").Value = 33.877666
$ws.Range("I5").Value = 0.5495662219753726
$ws.Range("J5").Value = 0.6375557499803809
$ws.Range("M5").Value = 19.878555
$ws.Range("N5").Value = 39.75711
$ws.Range("O5").Value = 0.1314800731212866
$ws.Range("P5").Value = 0.0916710195312133
$ws.Range("Q5").Value = 224.47968228421
$ws.Range("R5").Value = 1346.87809370526
$ws.Range("S5").Value = 0.0722570070503112
$ws.Range("T5").Value = 0.05844538560868884

$ws.Range("G6").Value = 11.29255533333333
$ws.Range("H6").Value = 33.877666
$ws.Range("I6").Value = 0.5495662219753726
$ws.Range("J6").Value = 0.6375557499803809
$ws.Range("M6").Value = 44.77944466666667
$ws.Range("N6").Value = 134.338334
$ws.Range("O6").Value = 0.2961787040911142
$ws.Range("P6").Value = 0.3097542059748472
$ws.Range("Q6").Value = 505.6743566942716
$ws.Range("R6").Value = 4551.069210248445
$ws.Range("S6").Value = 0.1627698114369155
$ws.Range("T6").Value = 0.1974855750998711

$ws.Range("I7").Value = 0.03478077306145753
$ws.Range("J7").Value = 0.04034942645199305
$ws.Range("M7").Value = 1.400501333333333
$ws.Range("N7").Value = 4.201504
$ws.Range("O7").Value = 0.00926314904242919
$ws.Range("P7").Value = 0.009687730200823723
$ws.Range("Q7").Value = 1.000909826072889
$ws.Range("R7").Value = 9.008188434656
$ws.Range("S7").Value = 0.0003221794846791873
$ws.Range("T7").Value = 0.0003908943572248887

$ws.Range("I8").Value = 0.03478077306145753
$ws.Range("J8").Value = 0.04034942645199305
$ws.Range("O8").Value = 0.1405812059498714
$ws.Range("P8").Value = 0.1470248171880475
$ws.Range("S8").Value = 0.004889523020848501
$ws.Range("T8").Value = 0.005932367047746846

$ws.Range("I9").Value = 0.03478077306145753
$ws.Range("J9").Value = 0.04034942645199305
$ws.Range("M9").Value = 63.87756733333333
$ws.Range("N9").Value = 191.632702
$ws.Range("O9").Value = 0.4224968677952986
$ws.Range("P9").Value = 0.4418622271050682
$ws.Range("Q9").Value = 45.65199852926423
$ws.Range("R9").Value = 410.867986763378
$ws.Range("S9").Value = 0.01469476767796491
$ws.Range("T9").Value = 0.0178288874344898

$ws.Range("I10").Value = 0.03478077306145753
$ws.Range("J10").Value = 0.04034942645199305
$ws.Range("M10").Value = 19.878555
$ws.Range("N10").Value = 39.75711
$ws.Range("O10").Value = 0.1314800731212866
$ws.Range("P10").Value = 0.0916710195312133
$ws.Range("Q10").Value = 14.206799061215
$ws.Range("R10").Value = 85.24079436729001
$ws.Range("S10").Value = 0.004572978585335311
$ws.Range("T10").Value = 0.00369887306035391

$ws.Range("I11").Value = 0.03478077306145753
$ws.Range("J11").Value = 0.04034942645199305
$ws.Range("M11").Value = 44.77944466666667
$ws.Range("N11").Value = 134.338334
$ws.Range("O11").Value = 0.2961787040911142
$ws.Range("P11").Value = 0.3097542059748472
$ws.Range("Q11").Value = 32.00295858789178
$ws.Range("R11").Value = 288.026627291026
$ws.Range("S11").Value = 0.01030132429262963
$ws.Range("T11").Value = 0.0124984045521776

$ws.Range("G12").Value = 8.5075845
$ws.Range("H12").Value = 17.015169
$ws.Range("I12").Value = 0.4140321595768645
$ws.Range("J12").Value = 0.3202144691088791
$ws.Range("M12").Value = 1.400501333333333
$ws.Range("N12").Value = 4.201504
$ws.Range("O12").Value = 0.00926314904242919
$ws.Range("P12").Value = 0.009687730200823723
$ws.Range("Q12").Value = 11.914883435696
$ws.Range("R12").Value = 71.489300614176
$ws.Range("S12").Value = 0.003835241602519322
$ws.Range("T12").Value = 0.003102151383126823

$ws.Range("G13").Value = 8.5075845
$ws.Range("H13").Value = 17.015169
$ws.Range("I13").Value = 0.4140321595768645
$ws.Range("J13").Value = 0.3202144691088791
$ws.Range("O13").Value = 0.1405812059498714
$ws.Range("P13").Value = 0.1470248171880475
$ws.Range("Q13").Value = 180.8249737179045
$ws.Range("R13").Value = 1084.949842307427
$ws.Range("S13").Value = 0.05820514029534521
$ws.Range("T13").Value = 0.04707947378170062

$ws.Range("G14").Value = 8.5075845
$ws.Range("H14").Value = 17.015169
$ws.Range("I14").Value = 0.4140321595768645
$ws.Range("J14").Value = 0.3202144691088791
$ws.Range("M14").Value = 63.87756733333333
$ws.Range("N14").Value = 191.632702
$ws.Range("O14").Value = 0.4224968677952986
$ws.Range("P14").Value = 0.4418622271050682
$ws.Range("Q14").Value = 543.443801742773
$ws.Range("R14").Value = 3260.662810456638
$ws.Range("S14").Value = 0.1749272905877485
$ws.Range("T14").Value = 0.1414906784717164

$ws.Range("G15").Value = 8.5075845
$ws.Range("H15").Value = 17.015169
$ws.Range("I15").Value = 0.4140321595768645
$ws.Range("J15").Value = 0.3202144691088791
$ws.Range("M15").Value = 19.878555
$ws.Range("N15").Value = 39.75711
$ws.Range("O15").Value = 0.1314800731212866
$ws.Range("P15").Value = 0.0916710195312133
$ws.Range("Q15").Value = 169.1184864003975
$ws.Range("R15").Value = 676.47394560159
$ws.Range("S15").Value = 0.05443697861573033
$ws.Range("T15").Value = 0.02935438685185716

$ws.Range("G16").Value = 8.5075845
$ws.Range("H16").Value = 17.015169
$ws.Range("I16").Value = 0.4140321595768645
$ws.Range("J16").Value = 0.3202144691088791
$ws.Range("M16").Value = 44.77944466666667
$ws.Range("N16").Value = 134.338334
$ws.Range("O16").Value = 0.2961787040911142
$ws.Range("P16").Value = 0.3097542059748472
$ws.Range("Q16").Value = 380.9649093647411
$ws.Range("R16").Value = 2285.789456188446
$ws.Range("S16").Value = 0.1226275084755211
$ws.Range("T16").Value = 0.09918777862047808

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.03330533333333333
$ws.Range("H17").Value = 0.099916
$ws.Range("I17").Value = 0.001620845386305282
$ws.Range("J17").Value = 0.001880354458746944
$ws.Range("M17").Value = 1.400501333333333
$ws.Range("N17").Value = 4.201504
$ws.Range("O17").Value = 0.00926314904242919
$ws.Range("P17").Value = 0.009687730200823723
$ws.Range("Q17").Value = 0.04664416374044444
$ws.Range("R17").Value = 0.419797473664
$ws.Range("S17").Value = 0.00001501413238807955
$ws.Range("T17").Value = 0.00001821636667825631

$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.03330533333333333
$ws.Range("H18").Value = 0.099916
$ws.Range("I18").Value = 0.001620845386305282
$ws.Range("J18").Value = 0.001880354458746944
$ws.Range("O18").Value = 0.1405812059498714
$ws.Range("P18").Value = 0.1470248171880475
$ws.Range("Q18").Value = 0.7078902389586667
$ws.Range("R18").Value = 6.371012150628
$ws.Range("S18").Value = 0.0002278603990650817
$ws.Range("T18").Value = 0.0002764587705459993

$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.03330533333333333
$ws.Range("H19").Value = 0.099916
$ws.Range("I19").Value = 0.001620845386305282
$ws.Range("J19").Value = 0.001880354458746944
$ws.Range("M19").Value = 63.87756733333333
$ws.Range("N19").Value = 191.632702
$ws.Range("O19").Value = 0.4224968677952986
$ws.Range("P19").Value = 0.4418622271050682
$ws.Range("Q19").Value = 2.127463672559111
$ws.Range("R19").Value = 19.147173053032
$ws.Range("S19").Value = 0.0006848020988944425
$ws.Range("T19").Value = 0.0008308576088888696

$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.03330533333333333
$ws.Range("H20").Value = 0.099916
$ws.Range("I20").Value = 0.001620845386305282
$ws.Range("J20").Value = 0.001880354458746944
$ws.Range("M20").Value = 19.878555
$ws.Range("N20").Value = 39.75711
$ws.Range("O20").Value = 0.1314800731212866
$ws.Range("P20").Value = 0.0916710195312133
$ws.Range("Q20").Value = 0.6620619004599999
$ws.Range("R20").Value = 3.97237140276
$ws.Range("S20").Value = 0.0002131088699097185
$ws.Range("T20").Value = 0.0001723740103133951

$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.03330533333333333
$ws.Range("H21").Value = 0.099916
$ws.Range("I21").Value = 0.001620845386305282
$ws.Range("J21").Value = 0.001880354458746944
$ws.Range("M21").Value = 44.77944466666667
$ws.Range("N21").Value = 134.338334
$ws.Range("O21").Value = 0.2961787040911142
$ws.Range("P21").Value = 0.3097542059748472
$ws.Range("Q21").Value = 1.491394331104889
$ws.Range("R21").Value = 13.422548979944
$ws.Range("S21").Value = 0.0004800598860479599
$ws.Range("T21").Value = 0.000582447702320423
